$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You manage state configurations of servers. Your company uses a cloud provider to host virtual servers, which are created and deleted hourly due to load changes. What method should you use to define an inventory file for this setup?",
        "ques_type": 2,
        "options": [
            "Dynamic inventory",
            "Liquid inventory",
            "Static inventory file and update it on each change",
            "The magic variable group"
        ],
        "score": "Dynamic inventory"
    },
    {
        "title": "You are responsible for a web server, webserver, that serves customers 24/7. You are writing a role that would update the configuration files of the webserver. Some updates of the configuration files require webserver restart. What method should you use to update the configuration files and restart webserver?",
        "ques_type": 2,
        "options": [
            "Add a webserver restart notify handler call for each task that updates any configuration file.",
            " Add a webserver restart task after each task that updates any configuration file.",
            "Add one webserver restart task at the end of the section that updates all web server configuration files.",
            "Manually restart webserver after each task that updates any configuration file."
        ],
        "score": "Add a webserver restart notify handler call for each task that updates any configuration file."
    },
    {
        "title": "You have a playbook update_config.yml. The playbook updates a configuration parameter string in a file at the desired server. However, this update is critical, and before applying the actual change, you want to see the changed string to make sure it updates correctly. Which command should you use?",
        "ques_type": 2,
        "options": [
            "ansible-playbook update_config.yml --list-tasks --check ",
            "ansible-playbook update_config.yml --diff --check",
            "ansible-playbook update_config.yml --list-tasks --syntax-check",
            "ansible-playbook update_config.yml --diff --list-tasks"
        ],
        "score": "ansible-playbook update_config.yml --diff --check"
    },
    {
        "title": "The variable non_existing_var is currently undefined.  True or false: The value of resulting_var will be equal to 1 after execution of the below tasks sequence. - set_fact:\n   intermediate_var: \"{{ non_existing_var | default(1) }}\n- set_fact:\n   resulting_var: \"{{ intermediate_var | default(2) }}\"",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    }
]
'@

# Clear any existing content/formatting on A1:A2, then set A1 to the new text
$ws.Range("A1:A2").Clear()
$ws.Range("A1").Value = $questionsText
